$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New athletes (column A) filled first for rows 4-7 ---
$ws.Range("A4").Value = "Иванов Иван Иванович"
$ws.Range("A5").Value = "Сидоров Петр Сергеевич"
$ws.Range("A6").Value = "Соловьев Андрей Петрович"
$ws.Range("A7").Value = "Сорокин Алексей Александрович"

# --- New teams (column H) filled next for rows 4-7 ---
$ws.Range("H4").Value = "МБОУ СОШ 6"
$ws.Range("H5").Value = "Лицей №7"
$ws.Range("H6").Value = "Гимназия ДГТУ"
$ws.Range("H7").Value = "СОШ №582"

# --- New mentors (column I) filled next for rows 4-7 ---
$ws.Range("I4").Value = "Антонов Алексей"
$ws.Range("I5").Value = "Смелов Егор"
$ws.Range("I6").Value = "Шмелев Вячислав"
$ws.Range("I7").Value = "Друзь Иван"

# --- Fill the remaining columns (birth date, payment, attendance, weight, group) ---
# Row 4
$ws.Range("B3").Copy($ws.Range("B4")) | Out-Null
$ws.Range("B4").Value = 38384
$ws.Range("C2:D2").Copy($ws.Range("C4")) | Out-Null
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = "A"

# Row 5
$ws.Range("B3").Copy($ws.Range("B5")) | Out-Null
$ws.Range("B5").Value = 38444
$ws.Range("C2:D2").Copy($ws.Range("C5")) | Out-Null
$ws.Range("E5").Value = 24
$ws.Range("F5").Value = "A"

# Row 6
$ws.Range("B3").Copy($ws.Range("B6")) | Out-Null
$ws.Range("B6").Value = 38902
$ws.Range("C2:D2").Copy($ws.Range("C6")) | Out-Null
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = "A"

# Row 7
$ws.Range("B3").Copy($ws.Range("B7")) | Out-Null
$ws.Range("B7").Value = 38058
$ws.Range("C2:D2").Copy($ws.Range("C7")) | Out-Null
$ws.Range("E7").Value = 26
$ws.Range("F7").Value = "A"

# --- Row 3 was missing the "Group" column; add it ---
$ws.Range("F3").Value = "A"

# --- Shorten the mentor name on row 2 (Таиров Дамир Шамилевич -> Таиров Дамир) ---
$ws.Range("I2").Value = "Таиров Дамир"

# --- View / zoom tweaks ---
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("B1").Select() | Out-Null

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
